$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8722577095031738
$ws.Range("B1").Value = 1.24306333065033
$ws.Range("C1").Value = 3.45703125
$ws.Range("D1").Value = 3.922667980194092
$ws.Range("E1").Value = 0.4346535205841064
